$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the repeated Q1-Q4 headers in the second and third question blocks
# (columns G:J and K:N) plus the final pair (O:P) with new, unique question
# labels Q5-Q14. This grows the shared string table with the new labels.
$ws.Range("G1").Value = "Q5"
$ws.Range("H1").Value = "Q6"
$ws.Range("I1").Value = "Q7"
$ws.Range("J1").Value = "Q8"
$ws.Range("K1").Value = "Q9"
$ws.Range("L1").Value = "Q10"
$ws.Range("M1").Value = "Q11"
$ws.Range("N1").Value = "Q12"
$ws.Range("O1").Value = "Q13"
$ws.Range("P1").Value = "Q14"

# Update the sheet's active selection to a single cell, Q3.
$ws.Range("Q3").Select()
